# Apply updates described in the commit:
# "Update AEP analysis in run_pipeline notebook and adjust results paths for
#  run_g1 and run_g3; modify metrics and outputs for enhanced multi-rule
#  simulations."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Text/content updates -------------------------------------------------

# Rule description for row 6: swap the swan p80 anomaly feature for waverys
$ws.Range("D6").Value = "swh_max_swan > t1 AND anom_swh_p80_waverys > t2"

# Threshold tuple for row 6
$ws.Range("H6").Value = "(np.float64(1.9166666), np.float64(0.1887734374999999))"

# Updated source_file names (new enhanced multi-rule "complete" CSV outputs)
$ws.Range("I6").Value = "enhanced_multi_rule_complete_summary_20250718_165152.csv"
$ws.Range("I7").Value = "enhanced_multi_rule_complete_summary_20250718_165101.csv"
$ws.Range("I9").Value = "enhanced_multi_rule_complete_summary_20250718_165039.csv"
$ws.Range("I10").Value = "enhanced_multi_rule_complete_summary_20250718_165124.csv"

# --- Metric updates: row 6 --------------------------------------------------

$ws.Range("E6").Value = 0.6532582461786002
$ws.Range("F6").Value = 119993.4537142857
$ws.Range("J6").Value = 406
$ws.Range("K6").Value = 330
$ws.Range("L6").Value = 1724
$ws.Range("M6").Value = 101
$ws.Range("N6").Value = 0.5516304347826086
$ws.Range("O6").Value = 0.8007889546351085
$ws.Range("P6").Value = 0.8317063647012886
$ws.Range("Q6").Value = 0.6532582461786002

# --- Metric updates: row 12 (TOTAL) ----------------------------------------

$ws.Range("E12").Value = 4.74831178606032
$ws.Range("F12").Value = 13286575.55294859
$ws.Range("J12").Value = 2040
$ws.Range("K12").Value = 2484
$ws.Range("L12").Value = 19835
$ws.Range("M12").Value = 1013
$ws.Range("N12").Value = 4.222917611953056
$ws.Range("O12").Value = 5.731080075981186
$ws.Range("P12").Value = 8.631240998323033
$ws.Range("Q12").Value = 4.74831178606032
